# From 1.2 to 1.3 version
# Swap the "Periodos Avaliativos" test-case block (rows 24-30, labeled TC3) with the
# "Avaliacoes" test-case block (rows 33-39, labeled TC4), so that the Avaliacoes
# content now appears first (as TC3) and the Periodos Avaliativos content follows (as TC4).
# The TC labels themselves (TC3 at B24, TC4 at B33) stay where they are; only the
# step/description text underneath is exchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Text currently describing "Periodos Avaliativos" (block starting at row 24 / TC3)
$periodosStep1   = $ws.Range("B28").Value2
$periodosTest1   = $ws.Range("D28").Value2
$periodosStep2   = $ws.Range("B29").Value2

# Text currently describing "Avaliacoes" (block starting at row 33 / TC4)
$avaliacoesStep1 = $ws.Range("B37").Value2
$avaliacoesTest1 = $ws.Range("D37").Value2
$avaliacoesStep2 = $ws.Range("B38").Value2

# Move the Avaliacoes content up into the TC3 block
$ws.Range("B28").Value2 = $avaliacoesStep1
$ws.Range("D28").Value2 = $avaliacoesTest1
$ws.Range("B29").Value2 = $avaliacoesStep2

# Move the Periodos Avaliativos content down into the TC4 block
$ws.Range("B37").Value2 = $periodosStep1
$ws.Range("D37").Value2 = $periodosTest1
$ws.Range("B38").Value2 = $periodosStep2
